$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink + content at D425 (it will be re-created at D475 below)
$ws.Range("D425").Hyperlinks.Delete()
$ws.Range("D425").Clear()

# Add the new translated rows (Sound, Floppy drive sound, Parallel/Serial/MIDI ports, etc.)
$ws.Range("A449").Value = 122
$ws.Range("B449").Value = 'Sound Emulation'
$ws.Range("C449").Value = 'Emulazione del suono'
$ws.Range("B450").Value = 'Disabled'
$ws.Range("C450").Value = 'Disabile'
$ws.Range("B451").Value = 'Disabled, but emulated'
$ws.Range("C451").Value = 'Disabilitato, ma emulato'
$ws.Range("B452").Value = 'Enabled'
$ws.Range("C452").Value = 'Abilitato'
$ws.Range("B453").Value = 'Automatic switching'
$ws.Range("C453").Value = 'Commutazione automatica'
$ws.Range("B454").Value = 'Include CD and FMV audio'
$ws.Range("C454").Value = 'Include CD e audio FMV'
$ws.Range("B455").Value = 'Volume Counter mode'
$ws.Range("C455").Value = 'Modalità del contatore di volume'
$ws.Range("B456").Value = 'Volume'
$ws.Range("C456").Value = 'Volume'
$ws.Range("B457").Value = 'Master'
$ws.Range("C457").Value = 'Master'
$ws.Range("B458").Value = 'Sound Buffer Size'
$ws.Range("C458").Value = 'Dimensione del buffer sonoro'
$ws.Range("B459").Value = 'Settings'
$ws.Range("C459").Value = 'Impostazioni'
$ws.Range("B460").Value = 'Channel mode:'
$ws.Range("C460").Value = 'Modalità di canale:'
$ws.Range("B461").Value = 'Stereo separation:'
$ws.Range("C461").Value = 'Separazione stereo:'
$ws.Range("B462").Value = 'Interpolation:'
$ws.Range("C462").Value = 'Interpolazione:'
$ws.Range("B463").Value = 'Frequency:'
$ws.Range("C463").Value = 'Frequenza:'
$ws.Range("B464").Value = 'Swap channels:'
$ws.Range("C464").Value = 'Scambiare i canali:'
$ws.Range("B465").Value = 'Stereo delay:'
$ws.Range("C465").Value = 'Ritardo stereo:'
$ws.Range("B466").Value = 'Audio filter:'
$ws.Range("C466").Value = 'Filtro audio:'
$ws.Range("B467").Value = 'Floppy Drive Sound Emulation'
$ws.Range("C467").Value = 'Emulazione sonora dell''unità floppy'
$ws.Range("B468").Value = 'Empty drive'
$ws.Range("C468").Value = 'Disco vuoto'
$ws.Range("B469").Value = 'Disk in drive'
$ws.Range("C469").Value = 'Disco nel drive'
$ws.Range("B470").Value = 'Drivers'
$ws.Range("C470").Value = 'Driver'
$ws.Range("B471").Value = 'DirectSound'
$ws.Range("C471").Value = 'DirectSound'
$ws.Range("B472").Value = 'WASAPI'
$ws.Range("C472").Value = 'WASAPI'
$ws.Range("B473").Value = 'OpenAL'
$ws.Range("C473").Value = 'OpenAL'
$ws.Range("B474").Value = 'PortAudio'
$ws.Range("C474").Value = 'PortAudio'
$ws.Range("A475").Value = 135
$ws.Range("B475").Value = 'Expand item'
$ws.Range("B476").Value = 'Leaf'
$ws.Range("B477").Value = 'Collapse Item'
$ws.Range("A478").Value = 138
$ws.Range("B478").Value = 'Parallel Port'
$ws.Range("C478").Value = 'Porta parallela'
$ws.Range("B479").Value = 'Printer:'
$ws.Range("C479").Value = 'Stampante'
$ws.Range("B480").Value = 'Type:'
$ws.Range("C480").Value = 'Tipo'
$ws.Range("B481").Value = 'Flush print job'
$ws.Range("C481").Value = 'Risciacquo del lavoro di stampa'
$ws.Range("B482").Value = 'Autoflush [] Time in seconds after a pending print job is automatically flushed.'
$ws.Range("C482").Value = 'Autoflush [] Tempo in secondi dopo il lavaggio automatico di un lavoro di stampa in attesa.'
$ws.Range("B483").Value = 'Ghostscript extra parameters:'
$ws.Range("C483").Value = 'Parametri extra di Ghostscript:'
$ws.Range("B484").Value = 'Sampler:'
$ws.Range("C484").Value = 'Campionario:'
$ws.Range("B485").Value = 'Stereo sampler'
$ws.Range("C485").Value = 'Campionatore stereo'
$ws.Range("B486").Value = 'Serial Port'
$ws.Range("C486").Value = 'Porta seriale'
$ws.Range("B487").Value = 'Shared'
$ws.Range("C487").Value = 'Condiviso'
$ws.Range("B488").Value = 'RTS/CTS'
$ws.Range("C488").Value = 'RTS/CTS'
$ws.Range("B489").Value = 'Direct []Use when emulating serial-link games on two PCs running WinUAE'
$ws.Range("C489").Value = 'Direct []Da utilizzare per l''emulazione di giochi con collegamento seriale su due PC con WinUAE.'
$ws.Range("B490").Value = 'uaeserial.device'
$ws.Range("C490").Value = 'uaeserial.device'
$ws.Range("B491").Value = 'MIDI'
$ws.Range("C491").Value = 'MIDI'
$ws.Range("B492").Value = 'Out'
$ws.Range("C492").Value = 'Fuori'
$ws.Range("B493").Value = 'IN'
$ws.Range("C493").Value = 'IN'
$ws.Range("B494").Value = 'Route MIDI In to MIDI Out'
$ws.Range("C494").Value = 'Indirizzare l''ingresso MIDI all''uscita MIDI'
$ws.Range("B495").Value = 'Protection Dongle'
$ws.Range("C495").Value = 'Dongle di protezione'
$ws.Range("A496").Value = 140
$ws.Range("B496").Value = 'UAE Authors and Contributors...'
$ws.Range("C496").Value = 'Autori e collaboratori degli Emirati Arabi Uniti...'
$ws.Range("A497").Value = 141
$ws.Range("B497").Value = 'Item1 - Item6'

# Re-create the DeepL hyperlink, now anchored at D475
$ws.Range("D475").Value = "DeepL Translate: The world's most accurate translator"
$ws.Hyperlinks.Add($ws.Range("D475"), "https://www.deepl.com/translator")
$ws.Range("D475").Style = "Hyperlink"

# Update the sheet selection / scroll position to match the new bottom of the data
$excel.Goto($ws.Range("A459"), $true)
$ws.Range("A498").Select()

